# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (only used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" colours   (used by the Slide Master / the deck's
#                                                   applied Design, "Integral")
#
# The authored change swaps the two themes' content: the deck-wide (Slide Master /
# Design) theme becomes the plain "Office Theme" colour set, while the Notes-Master-only
# theme keeps the "Integral" values that used to live on the Slide Master's theme part.
#
# The only theme part the PowerPoint object model lets us touch from script is the one
# backing the presentation's current Design (ppt/theme/theme2.xml here) - that's what
# ThemeColorScheme edits on a Slide/SlideRange resolve to. We drive all twelve theme
# colour slots there to the "Office Theme" palette so the applied-design theme matches
# the swapped-in colours from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as the BGR-packed OLE RGB values ThemeColorScheme.Item(n).RGB expects.
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
